# Weekly CompStat report refresh (111th Precinct) -- new crime data collected.
# Updates: report-header volume/date strings, one column width, and the
# weekly/28-day/YTD/2-year crime-statistics grid (rows 15-33).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1) Header text: "Volume 32   Number  19" -> "...Number  20"
#    and the reporting week "5/5/2025 ... 5/11/2025" -> "5/12/2025 ... 5/18/2025"
# ---------------------------------------------------------------------------
$ws.Range("A8").Value = "Volume 32   Number  20"
$ws.Range("C9").Value = "Report Covering the Week  5/12/2025  Through  5/18/2025"

# ---------------------------------------------------------------------------
# 2) Column H got a little wider to fit newly-longer percentage figures.
# ---------------------------------------------------------------------------
$ws.Columns.Item(8).ColumnWidth = 6.71

# ---------------------------------------------------------------------------
# 3) Crime-statistics grid updates.
# ---------------------------------------------------------------------------

# Row 15 - Rape
$ws.Range("L15").Value = -60

# Row 16 - Robbery
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = -60
$ws.Range("J16").Value = 23
$ws.Range("K16").Value = -60.869565217391
$ws.Range("L16").Value = -62.5
$ws.Range("M16").Value = -71.875
$ws.Range("N16").Value = -93.333333333333

# Row 17 - Fel. Assault
$ws.Range("D17").Value = 3
$ws.Range("F17").Value = 6
$ws.Range("G17").Value = 6
$ws.Range("J17").Value = 33
$ws.Range("K17").Value = 27.272727272727
$ws.Range("L17").Value = 13.513513513513
$ws.Range("N17").Value = -12.5

# Row 18 - Burglary
$ws.Range("C18").Value = 6
$ws.Range("D18").Value = 11
$ws.Range("E18").Value = -45.454545454545
$ws.Range("G18").Value = 21
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 109
$ws.Range("J18").Value = 101
$ws.Range("K18").Value = 7.920792079207
$ws.Range("L18").Value = -19.852941176470
$ws.Range("M18").Value = 12.371134020618
$ws.Range("N18").Value = -71.540469973890

# Row 19 - Gr. Larceny
$ws.Range("C19").Value = 7
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -22.222222222222
$ws.Range("F19").Value = 35
$ws.Range("G19").Value = 44
$ws.Range("H19").Value = -20.454545454545
$ws.Range("I19").Value = 158
$ws.Range("J19").Value = 185
$ws.Range("K19").Value = -14.594594594594
$ws.Range("L19").Value = -43.571428571428
$ws.Range("M19").Value = 10.489510489510
$ws.Range("N19").Value = -16.402116402116

# Row 20 - G.L.A.
$ws.Range("C20").Value = 6
$ws.Range("D20").Value = 10
$ws.Range("E20").Value = -40
$ws.Range("F20").Value = 25
$ws.Range("G20").Value = 23
$ws.Range("H20").Value = 8.695652173913
$ws.Range("I20").Value = 95
$ws.Range("J20").Value = 105
$ws.Range("K20").Value = -9.523809523809
$ws.Range("L20").Value = 69.642857142857
$ws.Range("M20").Value = 86.274509803921
$ws.Range("N20").Value = -92.135761589404

# Row 21 - TOTAL
$ws.Range("C21").Value = 19
$ws.Range("D21").Value = 35
$ws.Range("E21").Value = -45.714285714285
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 99
$ws.Range("H21").Value = -10.101010101010
$ws.Range("I21").Value = 417
$ws.Range("J21").Value = 450
$ws.Range("K21").Value = -7.333333333333
$ws.Range("L21").Value = -22.634508348794
$ws.Range("M21").Value = 20.520231213872
$ws.Range("N21").Value = -78.821736922295

# Row 24 - Petit Larceny
$ws.Range("C24").Value = 8
$ws.Range("D24").Value = 12
$ws.Range("E24").Value = -33.333333333333
$ws.Range("F24").Value = 37
$ws.Range("G24").Value = 43
$ws.Range("H24").Value = -13.953488372093
$ws.Range("I24").Value = 199
$ws.Range("J24").Value = 209
$ws.Range("K24").Value = -4.784688995215
$ws.Range("L24").Value = -2.450980392156
$ws.Range("M24").Value = 6.417112299465

# Row 25 - Retail Theft
# C25 switches from a numeric 1 to the literal text "0" (same style as the
# other "no data" cells in this table, e.g. A-column / C14). PasteSpecial is
# used so the destination keeps style index 13 (plain/general) instead of
# picking up a brand-new "quote-prefixed number" style that a plain
# Value="0" assignment would create.
$ws.Cells.Item(25, 1).Copy()                      # A25 already carries style 13
$ws.Cells.Item(25, 3).PasteSpecial(-4122)         # -4122 = xlPasteFormats
$ws.Cells.Item(14, 3).Copy()                      # C14 already holds the text "0"
$ws.Cells.Item(25, 3).PasteSpecial(-4163)         # -4163 = xlPasteValues (brings the "0" text)
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 8
$ws.Range("G25").Value = 7
$ws.Range("H25").Value = 14.285714285714
$ws.Range("J25").Value = 45
$ws.Range("K25").Value = -13.333333333333
$ws.Range("L25").Value = -13.333333333333

# Row 26 - Misd. Assault
$ws.Range("C26").Value = 3
$ws.Range("D26").Value = 8
$ws.Range("E26").Value = -62.5
$ws.Range("F26").Value = 18
$ws.Range("H26").Value = 0
$ws.Range("I26").Value = 67
$ws.Range("J26").Value = 76
$ws.Range("K26").Value = -11.842105263157
$ws.Range("L26").Value = -15.189873417721
$ws.Range("M26").Value = 17.543859649122

# Row 27 - UCR Rape*
$ws.Range("L27").Value = -33.333333333333

# Row 28 - Other Sex Crimes
# D28/G28 switch from text "0" to numeric 1 (style 13 -> 14); E28/H28 switch
# from text "***.*" to a numeric % figure (style 13 -> 15). Same two-step
# PasteSpecial trick: paste the target style first, then write the value.
$ws.Cells.Item(24, 4).Copy()                      # D24 already carries style 14
$ws.Cells.Item(28, 4).PasteSpecial(-4122)
$ws.Cells.Item(28, 4).Value = 1

$ws.Cells.Item(24, 5).Copy()                      # E24 already carries style 15
$ws.Cells.Item(28, 5).PasteSpecial(-4122)
$ws.Cells.Item(28, 5).Value = -100

$ws.Cells.Item(24, 7).Copy()                      # G24 already carries style 14
$ws.Cells.Item(28, 7).PasteSpecial(-4122)
$ws.Cells.Item(28, 7).Value = 1

$ws.Cells.Item(24, 8).Copy()                      # H24 already carries style 15
$ws.Cells.Item(28, 8).PasteSpecial(-4122)
$ws.Cells.Item(28, 8).Value = 100

$ws.Range("J28").Value = 4
$ws.Range("K28").Value = -50

# Row 29 - Shooting Vic.
# N29 switches from text "***.*" to numeric 100 (style 13 -> 15).
$ws.Cells.Item(24, 5).Copy()
$ws.Cells.Item(29, 14).PasteSpecial(-4122)
$ws.Cells.Item(29, 14).Value = 100

# Row 30 - Shooting Inc.
# N30 switches from text "***.*" to numeric 100 (style 13 -> 15).
$ws.Cells.Item(24, 5).Copy()
$ws.Cells.Item(30, 14).PasteSpecial(-4122)
$ws.Cells.Item(30, 14).Value = 100

# Row 33 - Traffic Fatalities
# C33/F33/I33 switch from text "0" to numeric 1 (style 13 -> 14).
$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item(33, 3).PasteSpecial(-4122)
$ws.Cells.Item(33, 3).Value = 1

$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item(33, 6).PasteSpecial(-4122)
$ws.Cells.Item(33, 6).Value = 1

$ws.Cells.Item(24, 4).Copy()
$ws.Cells.Item(33, 9).PasteSpecial(-4122)
$ws.Cells.Item(33, 9).Value = 1

$ws.Range("K33").Value = 0
$ws.Range("L33").Value = 0
